$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")
Write-Host ("I4 before = " + $ws.Range("I4").Value)
$ws.Range("I4").Value = 43
Write-Host ("I4 after = " + $ws.Range("I4").Value)
